$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 44.3284717741936
$ws.Range("D11").Value = 247.968

$ws.Range("C12").Value = 69.0645040322581
$ws.Range("D12").Value = 352.858

$ws.Range("C13").Value = 0.0620163690476191
$ws.Range("D13").Value = 10.509

$ws.Range("C15").Value = 56.872069892473
$ws.Range("D15").Value = 363.086

$ws.Range("C16").Value = 17.2449986559141
$ws.Range("D16").Value = 381.696
